# Krull_2003.xlsx -- "Updated soil type data"
#
# Adds a new controlled-vocabulary column (pro_usda_soil_order / USDA soil
# order) to the "profile" sheet, backed by a dropdown list sourced from a
# new column on the "controlled vocabulary" sheet, and populates the three
# data rows on "profile" with the soil order for each entry. Also wraps the
# long bibliographic-reference cell on the "metadata" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "controlled vocabulary" sheet: insert a new column E holding the
#    allowed USDA soil order values (the list that profile!N will pick
#    from). This pushes the former pro_soil_taxon_sys column (E) and
#    everything after it one column to the right.
# ---------------------------------------------------------------------
$cv = $wb.Worksheets.Item("controlled vocabulary")
$cv.Columns("E").Insert()

$cv.Range("E2").Value = "pro_usda_soil_order"

$soilOrders = @("Alfisols","Andisols","Aridisols","Entisols","Gelisols","Histosols","Inceptisols","Mollisols","Oxisols","Spodosols","Ultisols","Vertisols")
for ($i = 0; $i -lt $soilOrders.Length; $i++) {
    $cv.Cells.Item(4 + $i, 5).Value = $soilOrders[$i]
}

$cv.Range("E2:E15").Select()

# ---------------------------------------------------------------------
# 2. "profile" sheet: insert a new column N for pro_usda_soil_order
#    (shifting the former pro_soil_taxon column and everything after it
#    one column to the right), fill in the header + the three data rows,
#    and attach a dropdown list validation sourced from the controlled
#    vocabulary sheet.
# ---------------------------------------------------------------------
$pro = $wb.Worksheets.Item("profile")
$pro.Columns("N").Insert()

$pro.Range("N1").Value = "pro_usda_soil_order"
$pro.Range("N4").Value = "Oxisols"
$pro.Range("N5").Value = "Vertisols"
$pro.Range("N6").Value = "Vertisols"

$dvRange = $pro.Range("N4:N1048576")
$dvRange.Validation.Add(3, 1, 1, "='controlled vocabulary'!`$E`$4:`$E`$15")
$dvRange.Validation.IgnoreBlank = $true
$dvRange.Validation.InCellDropdown = $true
$dvRange.Validation.ShowInput = $true
$dvRange.Validation.ShowError = $true

$pro.Range("N6").Select()

# ---------------------------------------------------------------------
# 3. "metadata" sheet: wrap the long bibliographic reference text in M4
#    and grow row 4 to fit the wrapped text.
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("metadata")
$meta.Range("M4").WrapText = $true
$meta.Rows(4).RowHeight = 403.2

$meta.Activate()
$meta.Range("A4").Select()
